# Edit script: 
#  1) Re-style the three tables (on slides 14, 15, 16) that used table
#     style {2FEAE0F5-289D-4899-AB8D-4847CB07F1A4} so they use
#     {EA9633EE-AE1C-49F7-BF26-BD6569FBD691} instead.
#  2) Swap the presentation's colour theme: the slide master currently
#     carries the "Integral" (Red Violet) colour scheme while the default
#     "Office Theme" colours sit unused on the notes master. Re-apply the
#     Office Theme colours to the (visible) slide master so the deck's
#     on-screen theme becomes the plain Office palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style updates
# ---------------------------------------------------------------------
$oldStyle = "{2FEAE0F5-289D-4899-AB8D-4847CB07F1A4}"
$newStyle = "{EA9633EE-AE1C-49F7-BF26-BD6569FBD691}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colour swap: put the plain Office Theme colours onto the
#    slide master (the master that is actually used by every slide).
# ---------------------------------------------------------------------
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$master = $p.SlideMaster
$colorScheme = $master.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
